# Applies spell-checker "proofErr" run-splitting to several cells of the
# Arbeitsnachweis table, and appends a new row documenting 14.04.2023 work.
#
# Strategy: Word's InsertXML (as implemented here) replaces the ENTIRE
# paragraph that contains the target Range, so every replacement below
# supplies the complete, reconstructed paragraph content (all runs,
# including the ones that are not actually changing) rather than just the
# substring that differs.

$d = $word.ActiveDocument
$W = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-ParaText($p) {
    return $p.Range.Text.TrimEnd([char]7).TrimEnd([char]13)
}

function Set-ParaXml($p, [string]$innerXml) {
    $xml = '<w:p ' + $W + '>' + $innerXml + '</w:p>'
    $p.Range.InsertXML($xml)
}

# Walks $d.Paragraphs once, replacing each paragraph whose trimmed text
# exactly equals a key in $map with the corresponding inner-XML value.
# Matches are consumed left-to-right so repeated identical texts (e.g. the
# four standalone "Ament," paragraphs) are each handled independently and
# only once.
function Apply-Replacements($map) {
    $pending = @{}
    foreach ($k in $map.Keys) { $pending[$k] = New-Object System.Collections.Generic.Queue[string] }
    foreach ($k in $map.Keys) {
        foreach ($v in $map[$k]) { $pending[$k].Enqueue($v) }
    }

    $paras = @($d.Paragraphs)
    foreach ($p in $paras) {
        $t = Get-ParaText $p
        if ($pending.ContainsKey($t) -and $pending[$t].Count -gt 0) {
            $xml = $pending[$t].Dequeue()
            Set-ParaXml $p $xml
        }
    }
}

$replacements = [ordered]@{
    "Falk, Ament, Hollmann" = @(
        '<w:r><w:t xml:space="preserve">Falk, </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>Ament</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t>, Hollmann</w:t></w:r>'
    )

    "Issues erstellt (zu erledigende Aufgaben)" = @(
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>Issues</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> erstellt (zu erledigende Aufgaben)</w:t></w:r>'
    )

    "Klasse DBAPI_Base, DBAPI_Fahrplan und DBAPI_StaDa erstellt mit allen Funktionen" = @(
        '<w:r><w:t xml:space="preserve">Klasse </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>DBAPI_Base</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve">, </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>DBAPI_Fahrplan</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> und </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>DBAPI_StaDa</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> erstellt mit allen Funktionen</w:t></w:r>'
    )

    # Four distinct standalone paragraphs whose whole text is just "Ament,".
    "Ament," = @(
        ('<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>Ament</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t>,</w:t></w:r>'),
        ('<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>Ament</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t>,</w:t></w:r>'),
        ('<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>Ament</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t>,</w:t></w:r>'),
        ('<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>Ament</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t>,</w:t></w:r>')
    )

    "Bootstrap, jQuery und popper.js importiert" = @(
        '<w:r><w:t xml:space="preserve">Bootstrap, </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>jQuery</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> und popper.js </w:t></w:r>' +
        '<w:r><w:t>importiert</w:t></w:r>'
    )

    "Lesbarkeit verbessert, Auskommentierungen, Grundstruktur für Fahrplanauskunft und „Über uns“ eingefügt, Festinstallation der Dateien von jQuery, Bootstrap und popper.js;" = @(
        '<w:r><w:t>Lesbarkeit</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> verbessert</w:t></w:r>' +
        '<w:r><w:t>, Auskommentierungen, Grundstruktur für Fahrplanauskunft und „Über uns“</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> eingefügt</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">, </w:t></w:r>' +
        '<w:r><w:t>Festi</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">nstallation der Dateien von </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>jQuery</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t>, Bootstrap und popper.js</w:t></w:r>' +
        '<w:r><w:t>;</w:t></w:r>'
    )

    "Klasse MainHandler für alle API-Funktionen erstellt" = @(
        '<w:r><w:t xml:space="preserve">Klasse </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>MainHandler</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">für alle API-Funktionen </w:t></w:r>' +
        '<w:r><w:t>erstellt</w:t></w:r>'
    )

    "Request-Klassen erstellt, PanelBuilder erstellt" = @(
        '<w:r><w:t xml:space="preserve">Request-Klassen erstellt, </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>PanelBuilder</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> erstellt</w:t></w:r>'
    )

    'Nutzung von DevExpress-Komponenten („DevExtreme“) geplant' = @(
        '<w:r><w:t xml:space="preserve">Nutzung </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">von </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>DevExpress</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t>-Komponenten</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> („</w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>DevExtreme</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t>“)</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> geplant</w:t></w:r>'
    )

    "Implementierung der DevExpress-Komponenten," = @(
        '<w:r><w:t xml:space="preserve">Implementierung der </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>DevExpress</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t>-Komponenten,</w:t></w:r>'
    )

    "Download und Implementierung der JavaScript-DevExpress-Komponenten," = @(
        '<w:r><w:t>Download und Implementierung der JavaScript-</w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>DevExpress</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t>-Komponenten,</w:t></w:r>'
    )

    '„Paging“ beim Dropdown hinzugefügt (Daten werden nacheinander geladen)' = @(
        '<w:r><w:t>„</w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>Paging</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t>“ beim Dropdown hinzugefügt (Daten werden nacheinander geladen)</w:t></w:r>'
    )

    "Funktion erstellt, um eine Route zwischen zwei Bahnhöfen mit allen Zwischenstops von der API zu laden," = @(
        '<w:r><w:t xml:space="preserve">Funktion erstellt, um eine Route zwischen zwei Bahnhöfen mit allen </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>Zwischenstops</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> von der API zu laden,</w:t></w:r>'
    )
}

Apply-Replacements $replacements

# --- Append the new 14.04.2023 row -----------------------------------------

$t = $d.Tables.Item(1)
$newRow = $t.Rows.Add()

$cell1 = $newRow.Cells.Item(1)
$cell1.Range.InsertXML(
    '<w:p ' + $W + '><w:r><w:t>Kartenübersicht fortgesetzt,</w:t></w:r></w:p>' +
    '<w:p ' + $W + '><w:r><w:t>Begrenzung der Anfragen-Anzahl an die API</w:t></w:r></w:p>'
)

$cell2 = $newRow.Cells.Item(2)
$cell2.Range.InsertXML(
    '<w:p ' + $W + '><w:r><w:t>Hollmann</w:t></w:r></w:p>'
)

$cell3 = $newRow.Cells.Item(3)
$cell3.Range.InsertXML(
    '<w:p ' + $W + '><w:r><w:t>Anzahl der Anfragen an die API auf 100/Minute in der Query-Funktion begrenzt (Vorgabe der API),</w:t></w:r></w:p>' +
    '<w:p ' + $W + '>' +
        '<w:r><w:t xml:space="preserve">Kartenübersicht fortgesetzt: Anzeige der Weltkarte (Fokus auf Deutschland) über Bing mit Hilfe der </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>dxMap</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve">-Komponente von </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>DevExtreme</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '</w:p>'
)

$cell4 = $newRow.Cells.Item(4)
$cell4.Range.InsertXML(
    '<w:p ' + $W + '><w:r><w:t>1h</w:t></w:r></w:p>'
)

$cell5 = $newRow.Cells.Item(5)
$cell5.Range.InsertXML(
    '<w:p ' + $W + '><w:r><w:t>14.04.2023</w:t></w:r></w:p>'
)

Write-Host "Edit complete."
